# Fix: prevent hidden columns from being labeled upon detecting changes
# For rows where no real difference was detected in visible columns, clear the
# "Änderung" (L column) marker, and for the top-level "header" rows of each
# segment group, restore the unhighlighted (non-bold/grey) formatting across
# the whole row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Template row (row 2) already carries the desired target formatting:
#   columns A, C:K, M:V -> grey fill style
#   column B             -> bold + grey fill style
#   column L              -> grey fill, centered, no value
$templateRow = 2

$headerRows = @(123,126,130,132,136,141,149,155,159,163,166,169,175,182,186,189,192,195,200,203,205,208,211,214,218,221,224,231,237,240,243)
$detailRows = @(124,125,127,128,129,131,133,135,137,138,139,142,143,144,145,150,151,152,156,157,158,160,161,162,164,165,167,168,170,171,172,176,177,178,179,180,181,183,184,185,187,188,190,191,193,194,196,197,198,199,201,202,204,206,207,209,210,212,213,215,216,217,219,220,222,223,225,226,227,228,229,230,232,233,235,236,238,239,241,242,244)

foreach ($r in $headerRows) {
    $ws.Range("L$r").ClearContents()
    $ws.Range("A${templateRow}:V${templateRow}").Copy()
    $ws.Range("A${r}:V${r}").PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

foreach ($r in $detailRows) {
    $cell = $ws.Range("L$r")
    $cell.ClearContents()
    $ws.Range("L${templateRow}").Copy()
    $cell.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteFormats)
}

$excel.CutCopyMode = 0
